$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "30.593.25"
$c.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +0.56%  "

$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "1.920.82"
$c.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  -0.16%  "

$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = "1.0000"
$c.Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  -0.33%  "

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "247.90"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +3.47%  "

$ws.Cells.Item(6, 5).Value = "  -0.22%  "

$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "0.4745"
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +0.33%  "

$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "0.2891"
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +1.73%  "

$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "0.06835"
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +4.41%  "

$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "105.13"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +1.34%  "

$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "18.39"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -3.29%  "

$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "1.921.09"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -0.24%  "

$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "0.07700"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +1.60%  "

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "5.300"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +3.96%  "

$ws.Cells.Item(15, 5).Value = "  +3.18%  "

$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "291.45"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -0.38%  "

$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "30.597.23"
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +0.57%  "

$ws.Cells.Item(18, 5).Value = "  +1.55%  "

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "0.9997"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -0.25%  "

$ws.Cells.Item(20, 5).Value = "  +0.10%  "

$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "5.546"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +7.08%  "

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "2.168.61"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +0.10%  "

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "0.9999"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -0.36%  "

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "6.382"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +1.88%  "

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "9.407"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +1.99%  "

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "167.95"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +1.47%  "

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "21.11"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +8.78%  "

$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "2.112"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +4.34%  "

$ws.Cells.Item(29, 5).Value = "  -4.36%  "

$ws.Cells.Item(30, 5).Value = "  +2.85%  "

$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "4.179"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +2.15%  "

$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "4.071"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +4.18%  "

$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "0.05030"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +0.50%  "

$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "0.7386"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +0.62%  "

$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "1.144"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +0.10%  "

$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "0.02075"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +6.84%  "

$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "2.743"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +0.88%  "

$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "2.688"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -0.29%  "

$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "2.054"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +2.31%  "

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "111.04"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +3.82%  "

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "0.8778"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +0.77%  "

$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "0.4380"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +6.61%  "

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "5.873"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +0.92%  "

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "0.9997"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -0.26%  "

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "67.74"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -1.44%  "

$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "7.256"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +0.73%  "

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "9.296"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +0.65%  "

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "48.45"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +15.96%  "

$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "0.1233"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +2.69%  "

$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "34.80"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +1.00%  "

$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "0.2495"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +11.90%  "
